$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 4; $row -le 83; $row++) {
    $cell = $ws.Cells.Item($row, 8)   # Column H
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = $val * 0.322
    }
}
